$d = $word.ActiveDocument

# 1) "9/3" row, last cell: label text changes from the Lab2 entry to the Lab4 entry.
#    (The hyperlink itself keeps pointing at labs/Lab2-zoo/Lab2.qmd; only the
#    visible label + link text change, per the source diff.)
$d.Content.Find.Execute(
    "Lab2: PC2-G-Honolulu Zoo [", $true, $false, $false, $false, $false,
    $true, 1, $false, "Lab4: PC3-G-ECG/Cardiac Function [", 2) | Out-Null

# 2) "9/10" row, last cell: label text changes from the Lab3 entry to the Lab2 entry.
#    (The hyperlink keeps pointing at labs/Lab3-human-peripheral-circulation-dive-response/Lab3.qmd.)
$d.Content.Find.Execute(
    "Lab3:PC3-IWS-Peripheral Circ/Dive Resp [", $true, $false, $false, $false, $false,
    $true, 1, $false, "Lab2: PC2-G-Honolulu Zoo [", 2) | Out-Null

# 3) "9/17" row, last cell: label text changes from the Lab4 entry to the Lab3 entry.
#    (The hyperlink keeps pointing at labs/Lab4-human-ecg/Lab4.qmd; the visible
#    hyperlink text stays "overview" so no further edit is needed there.)
$d.Content.Find.Execute(
    "Lab Practical / Lab4: PC3-G-ECG/Cardiac Function [", $true, $false, $false, $false, $false,
    $true, 1, $false, "Lab Practical / Lab3:PC3-IWS-Peripheral Circ/Dive Resp [", 2) | Out-Null

# Now swap the two hyperlink display texts that moved with their labels:
#   labs/Lab2-zoo/Lab2.qmd        : "lab manual" -> "overview"
#   labs/Lab3-.../Lab3.qmd        : "overview"   -> "lab manual"
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.Address -eq "labs/Lab2-zoo/Lab2.qmd") {
        $h.TextToDisplay = "overview"
    } elseif ($h.Address -eq "labs/Lab3-human-peripheral-circulation-dive-response/Lab3.qmd") {
        $h.TextToDisplay = "lab manual"
    }
}

Write-Output "done"
